# Weekly refresh of the cryptocurrency price/volume snapshot (GitHub Actions job).
# Updates the "Price" (D) and "Volume(1h)" (E) columns for most rows, and for the
# Gas / MXToken pair the rows are also swapped back (B/C/D/E) to their new rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($a1, $text) {
    # Every cell in this sheet is stored as literal text (inline strings in the
    # source export), including "Price" values that happen to look numeric
    # (e.g. "246.33"). Assigning a numeric-looking string via .Value normally
    # makes Excel auto-convert the cell to a number, so we force Text format
    # first and clear the temporary formatting afterwards (ClearFormats resets
    # the style index back to the sheet default instead of leaving a stray
    # Text number-format applied to the cell).
    $cell = $ws.Range($a1)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# row -> updated columns (only the columns that actually changed)
$updates = @(
    @{ Row = 2;  D = "35.333.38";  E = "  +0.42%  " },
    @{ Row = 3;  D = "1.887.52";   E = "  -0.69%  " },
    @{ Row = 4;               E = "  -0.71%  " },
    @{ Row = 5;  D = "246.33";    E = "  -2.58%  " },
    @{ Row = 6;                E = "  -0.95%  " },
    @{ Row = 7;                E = "  -0.76%  " },
    @{ Row = 8;  D = "43.44";    E = "  +5.97%  " },
    @{ Row = 9;  D = "0.356";    E = "  -1.27%  " },
    @{ Row = 10; D = "54.06";    E = "  +2.17%  " },
    @{ Row = 11; D = "0.0741";   E = "  -1.62%  " },
    @{ Row = 12; D = "0.0972";   E = "  -1.24%  " },
    @{ Row = 13; D = "13.24";    E = "  +1.57%  " },
    @{ Row = 14; D = "2.159.97"; E = "  -0.74%  " },
    @{ Row = 15; D = "0.756";    E = "  +2.59%  " },
    @{ Row = 16; D = "1.880.60"; E = "  -0.86%  " },
    @{ Row = 17; D = "4.89";     E = "  -1.42%  " },
    @{ Row = 18; D = "35.396.51";E = "  +0.65%  " },
    @{ Row = 19; D = "73.03";    E = "  -0.92%  " },
    @{ Row = 20; D = "0.0₃0822"; E = "  -1.49%  " },
    @{ Row = 21; D = "244.46";   E = "  +0.84%  " },
    @{ Row = 22; D = "12.77";    E = "  -1.44%  " },
    @{ Row = 23; D = "4.94";     E = "  -2.16%  " },
    @{ Row = 24; D = "2.70";     E = "  +11.58%  " },
    @{ Row = 25;               E = "  -0.78%  " },
    @{ Row = 26; D = "2.14";     E = "  -6.07%  " },
    @{ Row = 27; D = "166.31";   E = "  -0.28%  " },
    @{ Row = 28; D = "8.50";     E = "  -1.10%  " },
    @{ Row = 29; D = "18.31";    E = "  -0.95%  " },
    @{ Row = 30;               E = "  -1.95%  " },
    @{ Row = 31; D = "4.128.44"; E = "  +0.00%  " },
    @{ Row = 32;               E = "  +10.33%  " },
    @{ Row = 33; D = "4.27";     E = "  -1.28%  " },
    @{ Row = 34; D = "0.0582";   E = "  -4.24%  " },
    @{ Row = 35; D = "4.17";     E = "  -0.97%  " },
    @{ Row = 36;               E = "  -0.80%  " },
    @{ Row = 37;               E = "  -12.28%  " },
    @{ Row = 38; D = "0.848";    E = "  -0.65%  " },
    @{ Row = 39; D = "1.96";     E = "  -2.37%  " },
    @{ Row = 40; D = "0.0694";   E = "  +6.88%  " },
    @{ Row = 41; D = "0.0220";   E = "  +2.78%  " },
    @{ Row = 42; D = "17.15";    E = "  -0.54%  " },
    @{ Row = 43; D = "96.53";    E = "  -4.65%  " },
    @{ Row = 44; D = "1.08";     E = "  -2.24%  " },
    @{ Row = 45; D = "1.299.84"; E = "  -1.14%  " },
    @{ Row = 46;               E = "  -4.53%  " },
    @{ Row = 47;               E = "  +7.42%  " },
    @{ Row = 48;               E = "  -1.08%  " },
    @{ Row = 49; B = "MXToken"; C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D = "2.73";  E = "  -0.72%  " },
    @{ Row = 50; B = "Gas";     C = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas";   D = "12.23"; E = "  +3.28%  " },
    @{ Row = 51;               E = "  -5.26%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B$row").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$row").Value = $u.C }
    if ($u.ContainsKey("D")) { Set-Text "D$row" $u.D }
    if ($u.ContainsKey("E")) { $ws.Range("E$row").Value = $u.E }
}
